# Apply scheduled-runner updates to the per-job Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each block below targets a specific job sheet/row and rewrites the computed
# market/profit columns (H..N) with freshly pulled pricing data.

$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 126.69231
$ws.Range("I15").Value = 126.69231
$ws.Range("K15").Value = 380.07693
$ws.Range("M15").Value = -211.07693

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1852.8889
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 1959.5
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 5878.5
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -7874.5

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1852.8889
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 1959.5
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 17635.5
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -27619.5

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7942.1665
$ws.Range("I86").Value = 7699
$ws.Range("K86").Value = 7699
$ws.Range("M86").Value = -6576

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 7942.1665
$ws.Range("I89").Value = 7699
$ws.Range("K89").Value = 38495
$ws.Range("M89").Value = -32879

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7872.5
$ws.Range("I138").Value = 6499.6665
$ws.Range("J138").Value = 8902.125
$ws.Range("K138").Value = 19498.9995
$ws.Range("L138").Value = 26706.375
$ws.Range("M138").Value = -14358.9995
$ws.Range("N138").Value = -36986.375

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5056.1924
$ws.Range("I32").Value = 5218.64
$ws.Range("J32").Value = 995
$ws.Range("K32").Value = 5218.64
$ws.Range("L32").Value = 995
$ws.Range("M32").Value = -4931.64
$ws.Range("N32").Value = -1569

# ARM row 51
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21512

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 967.8333
$ws.Range("I61").Value = 967.8333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 967.8333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -755.8333
$ws.Range("N61").ClearContents()

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 967.8333
$ws.Range("I136").Value = 967.8333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2903.4999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -353.4998999999998
$ws.Range("N136").ClearContents()

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1999
$ws.Range("I20").Value = 1999
$ws.Range("K20").Value = 1999
$ws.Range("M20").Value = -1752

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3378.7778
$ws.Range("I94").Value = 2081.8
$ws.Range("K94").Value = 2081.8
$ws.Range("M94").Value = -1630.8

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 100
$ws.Range("I134").Value = 100
$ws.Range("K134").Value = 300
$ws.Range("M134").Value = 2235

# CRP row 20
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 100000
$ws.Range("J20").Value = 100000
$ws.Range("L20").Value = 100000
$ws.Range("N20").Value = -100472

# CRP row 30
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 100000
$ws.Range("J30").Value = 100000
$ws.Range("L30").Value = 100000
$ws.Range("N30").Value = -100182

# CRP row 32
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 589.75
$ws.Range("I32").Value = 589.75
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 589.75
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -273.75
$ws.Range("N32").ClearContents()

# CRP row 55
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 40000
$ws.Range("I55").Value = 40000
$ws.Range("K55").Value = 40000
$ws.Range("M55").Value = -39685

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4184
$ws.Range("I58").Value = 1557.6
$ws.Range("J58").Value = 10750
$ws.Range("K58").Value = 1557.6
$ws.Range("L58").Value = 10750
$ws.Range("M58").Value = -1354.6
$ws.Range("N58").Value = -11156

# CRP row 70
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 24999.666
$ws.Range("J70").Value = 24999.666
$ws.Range("L70").Value = 24999.666
$ws.Range("N70").Value = -25629.666

# CRP row 73
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 24999.666
$ws.Range("J73").Value = 24999.666
$ws.Range("L73").Value = 24999.666
$ws.Range("N73").Value = -27183.666

# CRP row 128
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4184
$ws.Range("I136").Value = 1557.6
$ws.Range("J136").Value = 10750
$ws.Range("K136").Value = 4672.799999999999
$ws.Range("L136").Value = 32250
$ws.Range("M136").Value = -2122.799999999999
$ws.Range("N136").Value = -37350

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4499
$ws.Range("I133").Value = 4499
$ws.Range("K133").Value = 13497
$ws.Range("M133").Value = -8437

# GSM row 21
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 50000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# GSM row 30
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 50000
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

# LTW row 12
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 25500
$ws.Range("J12").Value = 25500
$ws.Range("L12").Value = 25500
$ws.Range("N12").Value = -25840

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6166
$ws.Range("I16").Value = 4500
$ws.Range("J16").Value = 6999
$ws.Range("K16").Value = 4500
$ws.Range("L16").Value = 6999
$ws.Range("M16").Value = -4330
$ws.Range("N16").Value = -7339

# LTW row 36
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# LTW row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 6943.8887
$ws.Range("J43").Value = 6943.8887
$ws.Range("L43").Value = 6943.8887
$ws.Range("N43").Value = -7329.8887

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1088.2157
$ws.Range("I46").Value = 250
$ws.Range("J46").Value = 4524.9
$ws.Range("K46").Value = 250
$ws.Range("L46").Value = 4524.9
$ws.Range("M46").Value = -62
$ws.Range("N46").Value = -4900.9

# LTW row 53
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1374.75
$ws.Range("I55").Value = 1250
$ws.Range("J55").Value = 1499.5
$ws.Range("K55").Value = 1250
$ws.Range("L55").Value = 1499.5
$ws.Range("M55").Value = -1077
$ws.Range("N55").Value = -1845.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1567.5
$ws.Range("I122").Value = 1567.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4702.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2252.5
$ws.Range("N122").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11110
$ws.Range("I132").Value = 11110
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 33330
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -30800
$ws.Range("N132").ClearContents()

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6939.7896
$ws.Range("I136").Value = 6873.5884
$ws.Range("K136").Value = 20620.7652
$ws.Range("M136").Value = -18070.7652

# WVR row 32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2199
$ws.Range("I132").Value = 1898.75
$ws.Range("J132").Value = 2599.3333
$ws.Range("K132").Value = 5696.25
$ws.Range("L132").Value = 7797.999899999999
$ws.Range("M132").Value = -3166.25
$ws.Range("N132").Value = -12857.9999

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1144
$ws.Range("I136").Value = 1148.7142
$ws.Range("K136").Value = 3446.1426
$ws.Range("M136").Value = -896.1425999999997

